# Update cryptos list prices and 1h volume percentages (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text representation (some prices contain
# "." thousands separators and must not be auto-converted to numbers by Excel).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.737.19"
$ws.Range("E2").Value = "  +3.57%  "

$ws.Range("D3").Value = "3.130.10"
$ws.Range("E3").Value = "  +2.19%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "590.11"
$ws.Range("E5").Value = "  +2.77%  "

$ws.Range("D6").Value = "146.24"
$ws.Range("E6").Value = "  +3.58%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "3.122.88"
$ws.Range("E8").Value = "  +2.07%  "

$ws.Range("E9").Value = "  +2.18%  "

$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +18.13%  "

$ws.Range("E11").Value = "  +4.02%  "

$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  +7.34%  "

$ws.Range("D14").Value = "36.11"
$ws.Range("E14").Value = "  +3.65%  "

$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").Value = "3.649.75"
$ws.Range("E16").Value = "  +2.26%  "

$ws.Range("D17").Value = "7.19"
$ws.Range("E17").Value = "  -0.35%  "

$ws.Range("D18").Value = "63.666.02"
$ws.Range("E18").Value = "  +3.56%  "

$ws.Range("D19").Value = "3.128.17"
$ws.Range("E19").Value = "  +2.14%  "

$ws.Range("D20").Value = "465.68"
$ws.Range("E20").Value = "  +3.90%  "

$ws.Range("D21").Value = "14.23"
$ws.Range("E21").Value = "  +2.43%  "

$ws.Range("D22").Value = "0.736"
$ws.Range("E22").Value = "  +1.09%  "

$ws.Range("D23").Value = "7.54"
$ws.Range("E23").Value = "  +2.08%  "

$ws.Range("D24").Value = "13.30"
$ws.Range("E24").Value = "  -1.28%  "

$ws.Range("D25").Value = "82.43"
$ws.Range("E25").Value = "  +0.93%  "

$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").Value = "8.72"
$ws.Range("E27").Value = "  +9.05%  "

$ws.Range("D28").Value = "2.71"
$ws.Range("E28").Value = "  +3.35%  "

$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  -0.02%  "

$ws.Range("D31").Value = "6.82"
$ws.Range("E31").Value = "  +2.68%  "

$ws.Range("E32").Value = "  +2.26%  "

$ws.Range("E33").Value = "  +0.65%  "

$ws.Range("D34").Value = "0.0₃0863"
$ws.Range("E34").Value = "  +9.23%  "

$ws.Range("E35").Value = "  +11.31%  "

$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +2.43%  "

$ws.Range("D37").Value = "3.39"
$ws.Range("E37").Value = "  +14.70%  "

$ws.Range("D38").Value = "6.13"
$ws.Range("E38").Value = "  +1.52%  "

$ws.Range("D39").Value = "50.91"
$ws.Range("E39").Value = "  +1.52%  "

$ws.Range("D40").Value = "449.83"
$ws.Range("E40").Value = "  +6.49%  "

$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("E42").Value = "  +1.71%  "

$ws.Range("D43").Value = "2.925.19"
$ws.Range("E43").Value = "  +5.74%  "

$ws.Range("E44").Value = "  +4.10%  "

$ws.Range("E45").Value = "  +2.79%  "

$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  +3.73%  "

$ws.Range("D47").Value = "126.36"
$ws.Range("E47").Value = "  +1.29%  "

$ws.Range("D49").Value = "34.90"
$ws.Range("E49").Value = "  -2.10%  "

$ws.Range("E50").Value = "  +0.55%  "

$ws.Range("D51").Value = "24.64"
$ws.Range("E51").Value = "  +3.63%  "

Write-Host "Updated cryptos list"